# Applies the "nach Planung in der Veranstaltung" update to the Story
# Cards workbook (Tabelle1): new acceptance-test notes, a couple of
# status/assignment corrections, a new final-acceptance row worth of
# data, and a refreshed cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 20: "Datenausgabe" story card gets a developer/acceptance note,
#     the acceptance-test text is refined and an estimated effort is added ---
$ws.Range("G20").Value = "Thomas, Mrosk/Schlufter"
$ws.Range("I20").Value = "Ergebnis speichern(xml oder txt)"
$ws.Range("K20").Value = "1h"

# --- Row 12/13: "Pseudocode"/"Ablaufplan" story cards get acceptance info ---
$ws.Range("I13").Value = "Ablaufplan und Formeln erstellen"
$ws.Range("J13").Value = "akzeptiert"

$ws.Range("I12").Value = "Ablaufplan studieren"
$ws.Range("J12").Value = "akzeptiert"

# --- Row 16: "Endabnahme" story card moves from jungfraeulich to in Arbeit,
#     gains a description / acceptance test / effort estimate / due date ---
$ws.Range("B16").Value = "in Arbeit"
$ws.Range("D16").Value = "Endabnahme"
$ws.Range("G16").Value = "alle"
$ws.Range("K16").Value = "-"
$ws.Range("L16").Value = "-"

$ws.Range("E16").Copy($ws.Range("M16"))
$ws.Range("M16").Value = 40847

# --- Row 17: "Öffentlichkeitsarbeit" story card moves to in Arbeit and
#     gets its description filled in ---
$ws.Range("B17").Value = "in Arbeit"
$ws.Range("D17").Value = "Öffentlichkeitsarbeit"

# --- refresh the remembered cursor position ---
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D24").Select()
